$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.264.98'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.592.50'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.47'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.08'
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0851'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = '1.817.28'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '1.596.46'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.88'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '26.240.42'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.05'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.35'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.97'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0491'
$ws.Range('E30').Value = '  -2.70%  '
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('D33').Value = '1.416.97'
$ws.Range('E33').Value = '  +6.76%  '
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.584'
$ws.Range('E37').Value = '  -3.24%  '
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.88'
$ws.Range('E40').Value = '  +2.91%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.975'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.765'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '1.729.05'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.99'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.46'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.49'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0954'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.13%  '
